$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.668.86"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "1.590.69"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.813.78"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").Value = "1.596.37"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "26.637.33"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "207.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.662"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +22.41%  "
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").Value = "1.322.53"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "1.726.71"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.836"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.02%  "
